$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save current values of row 7 (A7:E7)
$row7 = @()
for ($col = 1; $col -le 5; $col++) {
    $row7 += $ws.Cells.Item(7, $col).Value()
}

# Save current values of row 8 (A8:E8)
$row8 = @()
for ($col = 1; $col -le 5; $col++) {
    $row8 += $ws.Cells.Item(8, $col).Value()
}

# Swap: row 7 gets old row 8 values, row 8 gets old row 7 values
for ($col = 1; $col -le 5; $col++) {
    $ws.Cells.Item(7, $col).Value = $row8[$col - 1]
    $ws.Cells.Item(8, $col).Value = $row7[$col - 1]
}
